# Daily attendance processing - 2026-01-16 11:07:16
# Normalize the "Recorded By" column (G) so that "System" is listed
# before the recorder's email address, e.g.
#   "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
# Cells that don't contain "System" (e.g. just "dnasr281@gmail.com") are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the extent of the data so we only touch the column G cells
# that actually belong to the sheet's used range.
$lastRow = $ws.UsedRange.Rows.Count
$col = $ws.Range("G1:G$lastRow")

$changed = $col.Replace("dnasr281@gmail.com, System", "System, dnasr281@gmail.com")

Write-Host "Recorded-by order normalized in column G (changed: $changed)"
